# Update loading_percent values for Case_2_193 (380 kV case) — Sheet1
# Applies the 240 cell value changes (rows 2-25, columns B,C,E,F,G,H,I,K,N,O)
# as described by the commit "case with 380 kV done".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colMap = @{
    "B" = 2; "C" = 3; "E" = 5; "F" = 6; "G" = 7; "H" = 8; "I" = 9; "K" = 11; "N" = 14; "O" = 15
}

$data = @"
2,B,7.969501905455152
2,C,6.084422703784226
2,E,16.54411591767483
2,F,16.86991607391245
2,G,20.51507237767653
2,H,12.30069744046832
2,I,17.34948785135405
2,K,8.711901111011615
2,N,16.77778700103762
2,O,17.6064824593357
3,B,7.612625872626752
3,C,5.932719202731974
3,E,15.60481914985754
3,F,15.89584955866815
3,G,20.59319091330862
3,H,12.34732965362978
3,I,17.43962108630139
3,K,8.427911136948099
3,N,16.8230996316076
3,O,17.68480531990892
4,B,7.385332595964393
4,C,5.836752153496916
4,E,15.00272673997845
4,F,15.26997757108489
4,G,20.64924646487091
4,H,12.37787218702763
4,I,17.49829683317963
4,K,8.246984652499183
4,N,16.85245122500565
4,O,17.73684927412804
5,B,7.290778529108358
5,C,5.796970715052074
5,E,14.75125974661054
5,F,15.008197319934
5,G,20.67410983345193
5,H,12.39079900833029
5,I,17.52304620851023
5,K,8.1716751595724
5,N,16.86479775493279
5,O,17.75904955443883
6,B,7.27496529574907
6,C,5.790325398599843
6,E,14.70914371450483
6,F,14.96433081551589
6,G,20.67835996982557
6,H,12.39297452276874
6,I,17.52720647770467
6,K,8.159076741355706
6,N,16.86687119843541
6,O,17.76279574207554
7,B,7.384065045046853
7,C,5.836218328238304
7,E,14.999359712655
7,F,15.26647399323133
7,G,20.64957362016941
7,H,12.37804457695294
7,I,17.49862721634075
7,K,8.245975306744432
7,N,16.85261617245658
7,O,17.73714466139774
8,B,7.848215898451278
8,C,6.032721064041188
8,E,16.22565488309331
8,F,16.5399640634477
8,G,20.54032018964961
8,H,12.31637986973932
8,I,17.37987406823349
8,K,8.615380584365468
8,N,16.79309391290696
8,O,17.63266634848792
9,B,8.688976730472502
9,C,6.394263622634827
9,E,18.51991150758923
9,F,19.00274580682531
9,G,20.39086609009178
9,H,12.21060445998748
9,I,17.17343727168884
9,K,9.285130253034328
9,N,16.68846446862405
9,O,17.45925033676272
10,B,9.291192842729885
10,C,6.643654007435542
10,E,20.16011375280389
10,F,20.67494806633232
10,G,20.32132342866698
10,H,12.14211620115464
10,I,17.03786774348147
10,K,9.740800599173552
10,N,16.61890712498395
10,O,17.35115724205063
11,B,9.573210733458268
11,C,6.753265785869981
11,E,20.86377456740234
11,F,21.3917225636224
11,G,20.29857340001278
11,H,12.11296042591209
11,I,16.97968680761712
11,K,9.939662485897566
11,N,16.58883953747797
11,O,17.30620605337252
12,B,9.677646113287517
12,C,6.794197690229827
12,E,21.12416015724998
12,F,21.65686569030329
12,G,20.29124653252499
12,H,12.1022074228034
12,I,16.95815712335235
12,K,10.01371741606239
12,N,16.57767922034431
12,O,17.28979364064147
13,B,9.655259130647355
13,C,6.785408233898822
13,E,21.06835115681102
13,F,21.60004134736742
13,G,20.29276706345696
13,H,12.10451047940453
13,I,16.9627715868505
13,K,9.997824503372698
13,N,16.58007277166206
13,O,17.29330119152986
14,B,9.581850018512855
14,C,6.756644945241762
14,E,20.88531822358878
14,F,21.4136618050453
14,G,20.29794474371929
14,H,12.11207000430963
14,I,16.97790548147816
14,K,9.945780278603705
14,N,16.58791685395447
14,O,17.30484356117453
15,B,9.536577373750308
15,C,6.738950955617379
15,E,20.77241501713335
15,F,21.29868154950795
15,G,20.30128425713919
15,H,12.11673789551954
15,I,16.98724083573115
15,K,9.913737880073434
15,N,16.59275094137301
15,O,17.31199306651768
16,B,9.272432782746746
16,C,6.63641131097028
16,E,20.11327592698186
16,F,20.62722412089977
16,G,20.32298990342154
16,H,12.1440618483502
16,I,17.04174026035316
16,K,9.727631595594504
16,N,16.62090372719425
16,O,17.3541800710407
17,B,9.106198526225542
17,C,6.572506042820933
17,E,19.69805961939161
17,F,20.20408069597325
17,G,20.33858907771597
17,H,12.16133647299826
17,I,17.07606802039474
17,K,9.611274829016768
17,N,16.63857723694932
17,O,17.38114340088217
18,B,9.012600623841639
18,C,6.535389536220274
18,E,19.45523827114247
18,F,19.95656407809801
18,G,20.34839708041192
18,H,12.17146061225906
18,I,17.09614094585363
18,K,9.543558471972394
18,N,16.64889080705414
18,O,17.39704912627913
19,B,8.983742870179821
19,C,6.522761479357531
19,E,19.37233515454393
19,F,19.87204792380568
19,G,20.35186111186997
19,H,12.1749207989079
19,I,17.10299370748925
19,K,9.520496264180505
19,N,16.65240828789184
19,O,17.40250265818697
20,B,9.124053509714576
20,C,6.579346294948437
20,E,19.74267384416363
20,F,20.24955283636154
20,G,20.33684194119139
20,H,12.15947807789887
20,I,17.07237977083548
20,K,9.623743365858379
20,N,16.6366805255479
20,O,17.37823198748417
21,B,9.603476140458667
21,C,6.765109224256903
21,E,20.93924406267622
21,F,21.46857628470577
21,G,20.2963888982448
21,H,12.1098417809859
21,I,16.97344665848957
21,K,9.961101128087309
21,N,16.5856067394885
21,O,17.30143672103718
22,B,9.903055614602405
22,C,6.883150877927369
22,E,21.68587414653503
22,F,22.22866616901552
22,G,20.2774624093496
22,H,12.07907832149504
22,I,16.91171536899926
22,K,10.17428496795717
22,N,16.55354184059526
22,O,17.2548011881971
23,B,9.744425455613275
23,C,6.82046507291373
23,E,21.29061186321919
23,F,21.82633154458858
23,G,20.28687334874227
23,H,12.09534390760209
23,I,16.94439457488055
23,K,10.06118401950767
23,N,16.57053542937664
23,O,17.27936531992367
24,B,9.115986189770844
24,C,6.576254987935873
24,E,19.72251655142715
24,F,20.22900810905287
24,G,20.33762920635358
24,H,12.16031765797333
24,I,17.07404617735142
24,K,9.618108895616725
24,N,16.63753755262994
24,O,17.37954697875008
25,B,8.469545120792024
25,C,6.299192273888914
25,E,17.87801031873708
25,F,18.34778573295695
25,G,20.42427881740284
25,H,12.23759884959209
25,I,17.22645557613428
25,K,9.110137549535434
25,N,16.71548092939526
25,O,17.50278163412587
"@

$rows = $data -split "`n" | Where-Object { $_.Trim() -ne "" }
foreach ($line in $rows) {
    $parts = $line.Trim().Split(",")
    $rowNum = [int]$parts[0]
    $colLetter = $parts[1]
    $newVal = [double]$parts[2]
    $colNum = $colMap[$colLetter]
    $ws.Cells.Item($rowNum, $colNum).Value = $newVal
}
